# Append 23 days of weather observations (rows 366-388) to the Zanjan 1394 sheet,
# and restore the default view (select I7, no frozen/scrolled top-left cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shared lookup for the "weather condition" text column (Z).
$conditions = @("Rain-Thunderstorm", "Hail", "Rain", "Fog-Snow", "Thunderstorm", "Rain-Hail-Thunderstorm", "Snow", "Fog-Rain", "Fog", "Fog-Thunderstorm", "Rain-Snow", "Rain-Snow-Thunderstorm")

$newRows = @(
    @{Row=366; F=12; G=6; H=0; I=2; J=-3; K=-6; L=65; M=44; N=27; O=1015; P=1012; Q=1008; R=10; S=10; T=10; U=39; V=14; X=0.76; Y=6; AA=175},
    @{Row=367; F=14; G=9; H=5; I=4; J=3; K=2; L=87; M=64; N=39; O=1014; P=1011; Q=1007; R=10; S=9; T=6; U=26; V=5; X=2.03; Y=7; Z=0; AA=39},
    @{Row=368; F=13; G=7; H=0; I=2; J=-1; K=-4; L=87; M=53; N=18; O=1012; P=1010; Q=1006; R=10; S=10; T=8; U=26; V=8; X=0; Y=5; AA=324},
    @{Row=369; F=12; G=6; H=1; I=1; J=-4; K=-10; L=87; M=51; N=18; O=1019; P=1014; Q=1012; R=10; S=10; T=10; U=29; V=13; X=0.51; Y=5; AA=306},
    @{Row=370; F=17; G=9; H=0; I=1; J=-1; K=-4; L=85; M=44; N=23; O=1021; P=1019; Q=1013; R=10; S=10; T=10; U=23; V=8; X=0; Y=5; AA=177},
    @{Row=371; F=20; G=14; H=8; I=1; J=-3; K=-6; L=40; M=27; N=18; O=1020; P=1016; Q=1010; R=11; S=10; T=10; U=32; V=16; X=0; Y=7; AA=228},
    @{Row=372; F=13; G=11; H=9; I=4; J=1; K=-2; L=66; M=49; N=33; O=1018; P=1015; Q=1011; R=10; S=10; T=8; U=26; V=13; X=0; Y=8; AA=173},
    @{Row=373; F=10; G=8; H=7; I=7; J=6; K=5; L=93; M=82; N=70; O=1012; P=1009; Q=1007; R=10; S=7; T=1; U=35; V=13; X=7.11; Y=8; Z=2; AA=164},
    @{Row=374; F=12; G=8; H=3; I=5; J=3; K=1; L=85; M=69; N=44; O=1012; P=1009; Q=1007; R=10; S=8; T=6; U=26; V=10; X=9.91; Y=6; Z=2; AA=165},
    @{Row=375; F=11; G=6; H=1; I=1; J=-2; K=-6; L=88; M=54; N=29; O=1012; P=1010; Q=1007; R=10; S=9; T=2; U=39; V=16; X=7.11; Y=3; Z=6; AA=219},
    @{Row=376; F=11; G=6; H=1; I=0; J=-3; K=-7; L=87; M=50; N=27; O=1017; P=1014; Q=1010; R=10; S=9; T=6; U=26; V=10; X=0; Y=3; Z=6; AA=285},
    @{Row=377; F=11; G=4; H=-2; I=1; J=-2; K=-5; L=86; M=57; N=32; O=1018; P=1015; Q=1013; R=10; S=10; T=6; U=35; V=6; X=1.02; Y=5; Z=2; AA=239},
    @{Row=378; F=5; G=3; H=1; I=0; J=-4; K=-6; L=93; M=58; N=41; O=1024; P=1019; Q=1015; R=10; S=9; T=8; U=26; V=13; X=0; Y=7; Z=2; AA=36},
    @{Row=379; F=13; G=4; H=-5; I=0; J=-5; K=-10; L=86; M=47; N=21; O=1025; P=1020; Q=1016; R=14; S=11; T=10; U=23; V=8; X=0; Y=1; AA=328},
    @{Row=380; F=17; G=8; H=-2; I=-1; J=-5; K=-8; L=64; M=34; N=17; O=1017; P=1015; Q=1011; R=10; S=10; T=10; U=35; V=13; X=0; Y=5; AA=194},
    @{Row=381; F=14; G=9; H=5; I=2; J=1; K=-1; L=81; M=52; N=33; O=1015; P=1013; Q=1010; R=10; S=10; T=10; U=21; V=10; X=0; Y=5; AA=298},
    @{Row=382; F=19; G=10; H=1; I=2; J=-2; K=-4; L=75; M=40; N=21; O=1014; P=1013; Q=1009; R=10; S=10; T=10; U=23; V=6; X=0; Y=1; AA=235},
    @{Row=383; F=21; G=12; H=3; I=1; J=-2; K=-5; L=65; M=34; N=17; O=1016; P=1014; Q=1010; R=14; S=10; T=10; U=19; V=8; X=0; Y=1; AA=235},
    @{Row=384; F=21; G=11; H=2; I=2; J=-2; K=-6; L=70; M=34; N=16; O=1017; P=1015; Q=1011; R=10; S=10; T=8; U=35; V=10; X=0; Y=2; Z=2; AA=244},
    @{Row=385; F=18; G=11; H=4; I=11; J=3; K=0; L=96; M=59; N=28; O=1020; P=1017; Q=1014; R=10; S=9; T=7; U=32; V=6; X=0.51; Y=4; Z=2; AA=195},
    @{Row=386; F=20; G=11; H=2; I=4; J=1; K=-2; L=87; M=48; N=20; O=1019; P=1017; Q=1013; R=10; S=10; T=10; U=19; V=5; X=0; Y=4; AA=320},
    @{Row=387; F=19; G=12; H=5; I=8; J=3; K=1; L=83; M=53; N=29; O=1018; P=1016; Q=1013; R=11; S=9; T=7; U=35; V=13; X=4.06; Y=6; Z=2; AA=163},
    @{Row=388; F=17; G=12; H=7; I=9; J=8; K=5; L=88; M=74; N=43; O=1017; P=1015; Q=1009; R=10; S=9; T=6; U=26; V=10; X=1.02; Y=7; Z=2; AA=153}
)

$colIndex = @{
    A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13;
    N=14; O=15; P=16; Q=17; R=18; S=19; T=20; U=21; V=22; W=23; X=24; Y=25; Z=26; AA=27
}

foreach ($r in $newRows) {
    foreach ($col in $r.Keys) {
        if ($col -eq "Row") { continue }
        if ($col -eq "Z") {
            $ws.Cells.Item($r.Row, $colIndex[$col]).Value = $conditions[$r.Z]
        } else {
            $ws.Cells.Item($r.Row, $colIndex[$col]).Value = $r[$col]
        }
    }
}

# Scroll back to the top and select I7, matching the saved view after data entry.
$ws.Range("I7").Select() | Out-Null